# Update list of ANSPs for ACE: remove "HCAA" row from the ANSP sheet,
# which also shifts the remaining rows up (SMATSA's special bottom-border
# style on the old last row disappears along with the now-unused "HCAA"
# shared string). Also move the active tab from the data sheet to ANSP.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("ACE_landing_page_data")
$wsAnsp = $wb.Worksheets.Item("ANSP")

# The last row (SMATSA) had a unique bottom-border style reserved for the
# table's final row. Before removing a row above it, reformat it to match
# the regular ANSP-name style (same border as the rest of the list) so that
# once "HCAA" is removed, SMATSA - now the new last row - no longer carries
# a now-unused distinct style.
$lastCell = $wsAnsp.Range("A40")
$templateCell = $wsAnsp.Range("A39")
$templateCell.Copy()
$lastCell.PasteSpecial(-4122)

# Find and delete the row containing "HCAA" in column A.
$found = $wsAnsp.Range("A1:A40").Find("HCAA")
if ($found -ne $null) {
    $found.EntireRow.Delete()
}

$wsData.Range("A1:O7").Select()

# Select the new last data cell on the ANSP sheet and make it the active tab.
$wsAnsp.Activate()
$wsAnsp.Range("K15").Select()
